$d = $word.ActiveDocument
$d.Content.Find.Execute("Государственной итоговой атестации:", $true, $false, $false, $false, $false, $true, 1, $false, "Государственной итоговой аттестации:", 2)
